$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.677.99'
$ws.Range('E2').Value = '  -5.02%  '
$ws.Range('D3').Value = '2.185.82'
$ws.Range('E3').Value = '  -8.08%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '293.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '80.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.501'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.33%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.455'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -8.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0767'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -9.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '27.73'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.81'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -12.96%  '
$ws.Range('E13').Value = '  -2.43%  '
$ws.Range('D14').Value = '2.524.42'
$ws.Range('E14').Value = '  -8.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.08'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.77'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.29%  '
$ws.Range('D17').Value = '2.199.23'
$ws.Range('E17').Value = '  -8.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.701'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -8.11%  '
$ws.Range('D19').Value = '38.591.19'
$ws.Range('E19').Value = '  -5.06%  '
$ws.Range('E20').Value = '  -6.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '63.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -10.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '221.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.91%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -10.62%  '
$ws.Range('E27').Value = '  -5.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.48%  '
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '146.58'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '30.92'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -9.22%  '
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -10.31%  '
$ws.Range('E35').Value = '  -5.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0681'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.03%  '
$ws.Range('E37').Value = '  -5.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0938'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.11'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -12.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.63%  '
$ws.Range('D43').Value = '1.875.85'
$ws.Range('E43').Value = '  -4.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -16.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0253'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.48%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.84'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.42%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.64'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -11.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.52'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.61%  '
$ws.Range('D49').Value = '2.401.09'
$ws.Range('E49').Value = '  -7.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '69.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.12%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '85.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.14%  '
